# Replace "BROWN" with "H.BROWN" in column B (biosample harvester name),
# matching the shared-strings update in the diff, then update the saved
# selection to reflect column B being selected (B3:B27).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("B2:B27")
foreach ($cell in $rng.Cells) {
    if ($cell.Value() -eq "BROWN") {
        $cell.Value = "H.BROWN"
    }
}

$ws.Range("B3:B27").Select()
